# Rename the worksheet "Property1" -> "DataNode" (unifying DataNode / DataTable / Entity naming).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Move the active selection to D39, matching the saved cursor position.
$ws.Range("D39").Select() | Out-Null
